$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate a paragraph whose text (sans trailing paragraph mark) equals
# the given string. Searches from the end backwards so earlier edits (that
# may add/remove paragraphs later in the body) don't disturb indices we still
# need. Line breaks inside a paragraph show up as chr(11) in .Range.Text.
# ---------------------------------------------------------------------------
function Find-ParagraphByText {
    param([string]$targetText)
    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
        $p = $d.Paragraphs($i)
        $t = $p.Range.Text
        # Strip trailing control marks: paragraph mark (CR, 13) and, for the
        # last paragraph in a table cell, the cell-end mark (BEL, 7).
        while ($t.Length -gt 0 -and ([int][char]$t[$t.Length - 1] -eq 13 -or [int][char]$t[$t.Length - 1] -eq 7)) {
            $t = $t.Substring(0, $t.Length - 1)
        }
        if ($t -eq $targetText) {
            return $p
        }
    }
    return $null
}

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$w14Ns = "xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"
$NS = "$wNs $w14Ns"

# ---------------------------------------------------------------------------
# Change 3 (apply first, bottom-most edit): split the "If at least 1 FAIL,
# return to prev step " run so "prev" is wrapped in spell-check proofErr
# markers.
# ---------------------------------------------------------------------------
$lineBreak = [string][char]11
$target3 = "If at least 1 FAIL, return to prev step "
$p3 = Find-ParagraphByText $target3
if ($p3 -eq $null) { throw "Could not find paragraph for change 3" }
$full3 = $p3.Range

$xml3 = "<w:p $NS w14:paraId='4179FE90' w14:textId='7BD7A877' w:rsidR='00E42FFA' w:rsidRDefault='00E42FFA' w:rsidP='000754CA'>" +
        "<w:r><w:t xml:space='preserve'>If at least 1 FAIL, return to </w:t></w:r>" +
        "<w:proofErr w:type='spellStart'/>" +
        "<w:r><w:t>prev</w:t></w:r>" +
        "<w:proofErr w:type='spellEnd'/>" +
        "<w:r><w:t xml:space='preserve'> step </w:t></w:r>" +
        "</w:p>"
$full3.InsertXML($xml3) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: the "2. Test Applications" paragraph becomes an empty paragraph.
# ---------------------------------------------------------------------------
$target2 = "2. Test Applications"
$p2 = Find-ParagraphByText $target2
if ($p2 -eq $null) { throw "Could not find paragraph for change 2" }
$full2 = $p2.Range

$xml2 = "<w:p $NS></w:p>"
$full2.InsertXML($xml2) | Out-Null

# ---------------------------------------------------------------------------
# Change 1: split "Create configurations (create/delete/edit):" + linebreak +
# "1. Item types" into two paragraphs, adding an "Upload all Test
# Applications" line and renumbering "Item types" as step 2.
# ---------------------------------------------------------------------------
$target1 = "Create configurations (create/delete/edit):" + $lineBreak + "1. Item types"
$p1 = Find-ParagraphByText $target1
if ($p1 -eq $null) { throw "Could not find paragraph for change 1" }
$full1 = $p1.Range

$xml1 = "<w:p $NS w14:paraId='1ABAE9FA' w14:textId='3F210714' w:rsidR='00E42FFA' w:rsidRDefault='00E42FFA' w:rsidP='000754CA'>" +
        "<w:r><w:t>Create configurations (create/delete/edit):</w:t></w:r>" +
        "</w:p>" +
        "<w:p $NS>" +
        "<w:r><w:t>1</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'>. </w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'>Upload all </w:t></w:r>" +
        "<w:r><w:t>Test Applications</w:t></w:r>" +
        "<w:r><w:br/></w:r>" +
        "<w:r><w:t>2</w:t></w:r>" +
        "<w:r><w:t>. Item types</w:t></w:r>" +
        "</w:p>"
$full1.InsertXML($xml1) | Out-Null

Write-Output "done"
